# Updated cryptos list on Wed Jul 31 18:43:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to store the value as literal text (matching the
    # original inline-string cells) instead of letting Excel auto-convert
    # number-looking strings (e.g. "1.00", "26.66") into numeric values.
    # ClearFormats() afterwards drops the temporary "@" number format so
    # the cell's style index goes back to the sheet default.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.334.12"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.305.87"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "588.91"
$ws.Range("E5").Value = "  +2.69%  "

# Row 6 - Solana
Set-TextValue "D6" "179.56"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7 - XRP
Set-TextValue "D7" "0.641"
$ws.Range("E7").Value = "  +1.55%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.300.43"
$ws.Range("E9").Value = "  +0.37%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.125"
$ws.Range("E10").Value = "  -0.44%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.84"
$ws.Range("E11").Value = "  +2.08%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.04%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.875.40"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -2.47%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "66.354.38"
$ws.Range("E15").Value = "  +0.56%  "

# Row 16 - Avalanche
Set-TextValue "D16" "26.66"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17 - was WrappedEther, now ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000163"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18 - was ShibaInu, now WrappedEther
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.297.44"
$ws.Range("E18").Value = "  +0.68%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "425.92"
$ws.Range("E19").Value = "  -2.68%  "

# Row 20 - Polkadot
Set-TextValue "D20" "5.49"
$ws.Range("E20").Value = "  -1.79%  "

# Row 21 - Chainlink
Set-TextValue "D21" "13.03"
$ws.Range("E21").Value = "  -1.79%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.31"
$ws.Range("E22").Value = "  -1.42%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.01%  "

# Row 24 - Litecoin
Set-TextValue "D24" "71.17"
$ws.Range("E24").Value = "  -1.86%  "

# Row 25 - LEO
Set-TextValue "D25" "5.68"
$ws.Range("E25").Value = "  -0.19%  "

# Row 26 - Polygon
Set-TextValue "D26" "0.512"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +5.64%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +0.94%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "9.38"
$ws.Range("E29").Value = "  +5.01%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.02%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.75%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "22.31"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.06%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.18"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35 - Aptos
Set-TextValue "D35" "6.57"
$ws.Range("E35").Value = "  -0.93%  "

# Row 36 - Fetch.AI
Set-TextValue "D36" "1.19"
$ws.Range("E36").Value = "  -0.33%  "

# Row 37 - Monero
Set-TextValue "D37" "159.01"
$ws.Range("E37").Value = "  +0.14%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.24%  "

# Row 39 - Maker
$ws.Range("D39").Value = "2.859.99"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.80"
$ws.Range("E40").Value = "  +0.59%  "

# Row 41 - EnergySwap
Set-TextValue "D41" "26.19"
$ws.Range("E41").Value = "  -2.26%  "

# Row 42 - Filecoin
Set-TextValue "D42" "4.33"
$ws.Range("E42").Value = "  -0.20%  "

# Row 43 - Mantle
Set-TextValue "D43" "0.748"
$ws.Range("E43").Value = "  -4.19%  "

# Row 44 - OKB
Set-TextValue "D44" "39.61"

# Row 45 - RenderToken
Set-TextValue "D45" "5.89"
$ws.Range("E45").Value = "  -2.60%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  +0.48%  "

# Row 47 - Hedera
Set-TextValue "D47" "0.0638"
$ws.Range("E47").Value = "  -2.84%  "

# Row 48 - Bittensor
Set-TextValue "D48" "313.60"
$ws.Range("E48").Value = "  -2.44%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "22.90"
$ws.Range("E49").Value = "  -2.43%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +0.11%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.07%  "
